$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column J header + data -------------------------------------------
$ws.Range("J1").Value = "Projected Ranks"
# Give the new column a custom width similar to its neighbours (matches the
# sizing seen on the other data columns).
$ws.Columns.Item(10).ColumnWidth = 14.6

# --- Row 2 (NJ) gains the full set of metrics ------------------------------
$ws.Range("C2").Value = 7
$ws.Range("D2").Value = 579272
$ws.Range("E2").Value = 0.0679
$ws.Range("F2").Value = 470
$ws.Range("G2").Value = 11100000
$ws.Range("H2").Value = 0.36
$ws.Range("I2").Value = 2284
$ws.Range("J2").Value = 15

# --- Row 3 (OH) loses its detailed metrics, keeping only State + ITC ------
$ws.Range("C3:I3").Clear()

# --- Row 20 text changes; Row 21 is removed entirely -----------------------
$ws.Range("A20").Value = "** Solar by state seia.org"
$ws.Rows.Item(21).Delete()

# --- New threaded comments on the header cells that describe each metric --
$ws.Range("D1").AddCommentThreaded("State Homes Powered By Solar seia.org") | Out-Null
$ws.Range("F1").AddCommentThreaded("Solar Companies in State seia.org") | Out-Null
$ws.Range("G1").AddCommentThreaded("Total Solar Investment in State seia.org") | Out-Null
$ws.Range("H1").AddCommentThreaded("Prices have Fallen seia.org") | Out-Null
$ws.Range("I1").AddCommentThreaded("Growth Projections in the next 5 years seia.org") | Out-Null
$ws.Range("J1").AddCommentThreaded("Growth Projection Rankings  in the next 5 Years seia.org") | Out-Null

Write-Host "Edit complete"
